$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Date column to be treated as text so "2025-06-18" is not
# auto-converted into a date serial number (matches the rest of column E).
$ws.Range("E16").NumberFormat = "@"

$ws.Range("A16").Value = "edit1"
$ws.Range("B16").Value = "riya-morankar"
$ws.Range("C16").Value = "Merged"
$ws.Range("D16").Value = "N/A"
$ws.Range("E16").Value = "2025-06-18"
$ws.Range("F16").Value = "15431bee5d95560e049c0d32a92e0cea477e3ba5"
